# This script updates the "想去人数" (interest count) column F
# across all four worksheets to reflect refreshed scrape data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3725
$ws.Range("F5").Value = 3725
$ws.Range("F6").Value = 290
$ws.Range("F7").Value = 5261
$ws.Range("F8").Value = 576
$ws.Range("F9").Value = 408
$ws.Range("F11").Value = 1036
$ws.Range("F13").Value = 127
$ws.Range("F14").Value = 45
$ws.Range("F15").Value = 723
$ws.Range("F16").Value = 349
$ws.Range("F17").Value = 43
$ws.Range("F19").Value = 169
$ws.Range("F21").Value = 369
$ws.Range("F22").Value = 6017
$ws.Range("F26").Value = 6321
$ws.Range("F29").Value = 3250
$ws.Range("F31").Value = 742
$ws.Range("F32").Value = 4454
$ws.Range("F34").Value = 133
$ws.Range("F35").Value = 148
$ws.Range("F36").Value = 1115
$ws.Range("F40").Value = 911
$ws.Range("F41").Value = 1105
$ws.Range("F42").Value = 2052

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 50
$ws.Range("F5").Value = 63

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1146

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1146
$ws.Range("F7").Value = 3725
$ws.Range("F8").Value = 3725
$ws.Range("F9").Value = 290
$ws.Range("F10").Value = 5261
$ws.Range("F11").Value = 576
$ws.Range("F12").Value = 408
$ws.Range("F14").Value = 1036
$ws.Range("F16").Value = 127
$ws.Range("F17").Value = 45
$ws.Range("F18").Value = 723
$ws.Range("F19").Value = 349
$ws.Range("F20").Value = 43
$ws.Range("F21").Value = 50
$ws.Range("F23").Value = 169
$ws.Range("F25").Value = 369
$ws.Range("F26").Value = 6017
$ws.Range("F30").Value = 6321
$ws.Range("F33").Value = 3250
$ws.Range("F35").Value = 742
$ws.Range("F36").Value = 4454
$ws.Range("F39").Value = 133
$ws.Range("F40").Value = 148
$ws.Range("F41").Value = 1115
$ws.Range("F45").Value = 911
$ws.Range("F46").Value = 1105
$ws.Range("F48").Value = 2052
$ws.Range("F50").Value = 63
